# Append a new row 41 to each of the four worksheets, duplicating row 40
# (columns B..I identical) but with an updated timestamp in column A.
$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{ Sheet = 1; A = 45744.82736297454; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0xd"; G = [double]"5.68631262647114e+23"; I = 13 },
    @{ Sheet = 2; A = 45744.67373422454; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0xe"; G = [double]"5.68631262647114e+23"; I = 14 },
    @{ Sheet = 3; A = 45744.81831989584; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x3"; G = [double]"5.68631262647114e+23"; I = 3 },
    @{ Sheet = 4; A = 45744.87477822917; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; E = "0x3"; G = [double]"9.85046333984776e+23"; I = 3 }
)

foreach ($item in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($item.Sheet)

    $ws.Cells.Item(41, 1).Value = $item.A
    $ws.Cells.Item(41, 2).Value = "0x01,0x90"
    $ws.Cells.Item(41, 3).Value = $item.C
    $ws.Cells.Item(41, 4).Value = "0x01,0x7a"
    $ws.Cells.Item(41, 5).Value = $item.E
    $ws.Cells.Item(41, 6).Value = 400
    $ws.Cells.Item(41, 7).Value = $item.G
    $ws.Cells.Item(41, 8).Value = 378
    $ws.Cells.Item(41, 9).Value = $item.I

    # Match the number format used by the other rows in column A (datetime style).
    $ws.Cells.Item(41, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
